$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''42.484.39'
$ws.Range("E2").Value = '  -1.42%  '

$ws.Range("D3").Value = '''2.186.76'
$ws.Range("E3").Value = '  -2.11%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = '''250.77'
$ws.Range("E5").Value = '  +1.74%  '

$ws.Range("D6").Value = '''0.612'
$ws.Range("E6").Value = '  -1.12%  '

$ws.Range("D7").Value = '''74.80'
$ws.Range("E7").Value = '  -0.78%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  -5.19%  '

$ws.Range("D10").Value = '''40.25'
$ws.Range("E10").Value = '  -1.32%  '

$ws.Range("D11").Value = '''0.0908'
$ws.Range("E11").Value = '  -2.34%  '

$ws.Range("D12").Value = '''0.102'
$ws.Range("E12").Value = '  +0.28%  '

$ws.Range("D13").Value = '''6.77'
$ws.Range("E13").Value = '  -2.65%  '

$ws.Range("D14").Value = '''2.515.87'
$ws.Range("E14").Value = '  -2.15%  '

$ws.Range("D15").Value = '''14.18'
$ws.Range("E15").Value = '  -3.83%  '

$ws.Range("D16").Value = '''2.191.48'
$ws.Range("E16").Value = '  -0.96%  '

$ws.Range("D17").Value = '''0.768'
$ws.Range("E17").Value = '  -5.31%  '

$ws.Range("D18").Value = '''42.416.04'
$ws.Range("E18").Value = '  -1.38%  '

$ws.Range("E19").Value = '  -3.20%  '

$ws.Range("D20").Value = '''70.79'

$ws.Range("E21").Value = '  -2.17%  '

$ws.Range("D22").Value = '''226.95'
$ws.Range("E22").Value = '  -1.51%  '

$ws.Range("D23").Value = '''9.44'
$ws.Range("E23").Value = '  -10.12%  '

$ws.Range("D24").Value = '''2.12'
$ws.Range("E24").Value = '  -2.73%  '

$ws.Range("E26").Value = '  -4.52%  '

$ws.Range("D27").Value = '''3.38'
$ws.Range("E27").Value = '  +0.60%  '

$ws.Range("D28").Value = '''2.20'
$ws.Range("E28").Value = '  +0.10%  '

$ws.Range("D29").Value = '''2.15'
$ws.Range("E29").Value = '  -4.42%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = '''37.40'
$ws.Range("E30").Value = '  +0.94%  '

$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").Value = '''172.46'
$ws.Range("E31").Value = '  -1.26%  '

$ws.Range("E32").Value = '  -1.58%  '

$ws.Range("D33").Value = '''0.0812'
$ws.Range("E33").Value = '  +2.28%  '

$ws.Range("D34").Value = '''5.13'
$ws.Range("E34").Value = '  -4.37%  '

$ws.Range("E35").Value = '  -1.60%  '

$ws.Range("D36").Value = '''0.107'
$ws.Range("E36").Value = '  -3.41%  '

$ws.Range("E37").Value = '  -3.32%  '

$ws.Range("D38").Value = '''0.0333'
$ws.Range("E38").Value = '  +0.23%  '

$ws.Range("D39").Value = '''12.02'
$ws.Range("E39").Value = '  -8.30%  '

$ws.Range("E40").Value = '  -3.29%  '

$ws.Range("D41").Value = '''2.59'
$ws.Range("E41").Value = '  +12.34%  '

$ws.Range("E42").Value = '  -7.65%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = '''0.193'
$ws.Range("E43").Value = '  -2.69%  '

$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D44").Value = '''58.57'
$ws.Range("E44").Value = '  -2.56%  '

$ws.Range("D45").Value = '''101.42'
$ws.Range("E45").Value = '  -3.97%  '

$ws.Range("D46").Value = '''0.0972'
$ws.Range("E46").Value = '  -1.97%  '

$ws.Range("B47").Value = 'WOONetwork'
$ws.Range("C47").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D47").Value = '''0.460'
$ws.Range("E47").Value = '  +4.05%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").Value = '''8.16'
$ws.Range("E48").Value = '  -4.60%  '

$ws.Range("E49").Value = '  -1.74%  '

$ws.Range("D50").Value = '''1.13'
$ws.Range("E50").Value = '  -2.41%  '

$ws.Range("E51").Value = '  -0.69%  '
